$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

# Fix sorting / update average time values for rows with Rows=5000 and Rows=10000
$ws.Range("D2").Value = 21.676412
$ws.Range("D3").Value = 86.153375

$wb.RefreshAll()
$excel.CalculateFullRebuild()
